$d = $word.ActiveDocument
$r = $d.Content

# The target paragraph is the unique list item about following free online
# developer-training courses; locate it by its (unique) leading text.
$found = $r.Find.Execute("Suivre des formations gratuites en ligne sur le", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $para = $r.Paragraphs(1).Range

    $targetXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="5DAF53AA" w14:textId="77777777" w:rsidR="001408ED" w:rsidRPr="004D4F7D" w:rsidRDefault="001408ED" w:rsidP="001408ED"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:line="276" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/></w:rPr><w:t xml:space="preserve">Suivre des formations gratuites en ligne sur le </w:t></w:r><w:r w:rsidR="004D4F7D"><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/></w:rPr><w:t>d&#233;veloppement</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/></w:rPr><w:t xml:space="preserve"> informatique</w:t></w:r><w:r w:rsidR="00141F7C"><w:rPr><w:rFonts w:ascii="Century Gothic" w:hAnsi="Century Gothic"/></w:rPr><w:t>,</w:t></w:r></w:p>'

    $para.InsertXML($targetXml)
}
